# Updates cryptos list cell values (Price / Volume(1h) columns, plus the
# Bittensor/Kaspa row swap) to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need a leading
# apostrophe so Excel stores them as text (matching the source file's
# inline-string cells) instead of silently converting to a Double and
# dropping significant trailing/leading zeros.

$ws.Range("D2").Value = '60.448.33'
$ws.Range("E2").Value = '  -4.19%  '
$ws.Range("D3").Value = '2.965.34'
$ws.Range("E3").Value = '  -6.65%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''557.19'
$ws.Range("E5").Value = '  -5.89%  '
$ws.Range("D6").Value = '''123.95'
$ws.Range("E6").Value = '  -8.70%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '2.962.45'
$ws.Range("E8").Value = '  -6.66%  '
$ws.Range("D9").Value = '''0.491'
$ws.Range("E9").Value = '  -4.49%  '
$ws.Range("D10").Value = '''0.132'
$ws.Range("E10").Value = '  -7.46%  '
$ws.Range("D11").Value = '''5.09'
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("D12").Value = '''0.432'
$ws.Range("E12").Value = '  -5.13%  '
$ws.Range("D13").Value = '''0.0000220'
$ws.Range("E13").Value = '  -7.77%  '
$ws.Range("D14").Value = '''32.20'
$ws.Range("E14").Value = '  -7.37%  '
$ws.Range("D15").Value = '''0.118'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = '3.455.13'
$ws.Range("E16").Value = '  -6.61%  '
$ws.Range("D17").Value = '60.585.55'
$ws.Range("E17").Value = '  -3.93%  '
$ws.Range("D18").Value = '2.964.03'
$ws.Range("E18").Value = '  -6.73%  '
$ws.Range("D19").Value = '''6.08'
$ws.Range("E19").Value = '  -7.82%  '
$ws.Range("D20").Value = '''426.49'
$ws.Range("E20").Value = '  -7.72%  '
$ws.Range("D21").Value = '''12.91'
$ws.Range("E21").Value = '  -7.41%  '
$ws.Range("D22").Value = '''0.653'
$ws.Range("E22").Value = '  -7.26%  '
$ws.Range("D23").Value = '''7.06'
$ws.Range("E23").Value = '  -7.64%  '
$ws.Range("D24").Value = '''12.82'
$ws.Range("E24").Value = '  -4.40%  '
$ws.Range("D25").Value = '''77.89'
$ws.Range("E25").Value = '  -6.84%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '''2.45'
$ws.Range("E28").Value = '  -8.71%  '
$ws.Range("D29").Value = '''7.04'
$ws.Range("E29").Value = '  -9.39%  '
$ws.Range("D30").Value = '''1.86'
$ws.Range("E30").Value = '  -8.45%  '
$ws.Range("D31").Value = '''25.02'
$ws.Range("E31").Value = '  -8.24%  '
$ws.Range("D32").Value = '''5.92'
$ws.Range("E32").Value = '  -12.59%  '
$ws.Range("D33").Value = '''0.0922'
$ws.Range("E33").Value = '  -11.47%  '
$ws.Range("D34").Value = '''2.24'
$ws.Range("E34").Value = '  -6.07%  '
$ws.Range("D35").Value = '''0.944'
$ws.Range("E35").Value = '  -9.39%  '
$ws.Range("D36").Value = '''5.46'
$ws.Range("E36").Value = '  -6.33%  '
$ws.Range("D37").Value = '''49.32'
$ws.Range("E37").Value = '  -3.96%  '
$ws.Range("D38").Value = '0.0₃0654'
$ws.Range("E38").Value = '  -7.57%  '
$ws.Range("D39").Value = '''0.0353'
$ws.Range("E39").Value = '  -9.08%  '
$ws.Range("D40").Value = '''7.69'
$ws.Range("E40").Value = '  -5.40%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.106'
$ws.Range("E41").Value = '  -6.34%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '''370.33'
$ws.Range("E42").Value = '  -8.79%  '
$ws.Range("D43").Value = '2.653.22'
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("D44").Value = '''2.41'
$ws.Range("E44").Value = '  -8.58%  '
$ws.Range("D46").Value = '''0.231'
$ws.Range("E46").Value = '  -8.36%  '
$ws.Range("D47").Value = '''118.56'
$ws.Range("E47").Value = '  -4.51%  '
$ws.Range("D48").Value = '''32.68'
$ws.Range("E48").Value = '  -4.85%  '
$ws.Range("D49").Value = '''1.93'
$ws.Range("E49").Value = '  -9.23%  '
$ws.Range("D50").Value = '''0.105'
$ws.Range("E50").Value = '  -5.81%  '
$ws.Range("D51").Value = '''23.03'
$ws.Range("E51").Value = '  -9.82%  '
